# Add a new "NewContact" worksheet after the existing "Contacts" sheet,
# populate it with a new contact row (plus header row), and make it the
# active sheet - mirroring a manual "Add files via upload" style edit
# that introduced a second sheet to the FreeCRM workbook.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet immediately after the last existing sheet so it
# lands at the end of the tab strip (after "Contacts").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "NewContact"

# Populate the data row first (row 2) so the shared-string table picks up
# the new contact's values before the new header labels, then fill in the
# header row (row 1) - reusing the existing "Company"/"Description"/
# "Position"/"Department" headers from the Contacts sheet.
$newSheet.Range("A2").Value = "Nagendra"
$newSheet.Range("B2").Value = "Prasad"
$newSheet.Range("C2").Value = "Private Org"
$newSheet.Range("D2").Value = "It is a Test"
$newSheet.Range("E2").Value = "Engineer"
$newSheet.Range("F2").Value = "Networking"

$newSheet.Range("A1").Value = "FirstName"
$newSheet.Range("B1").Value = "LastName"
$newSheet.Range("C1").Value = "Company"
$newSheet.Range("D1").Value = "Description"
$newSheet.Range("E1").Value = "Position"
$newSheet.Range("F1").Value = "Department"

# Auto-fit the columns to the new content (as Excel does for a freshly
# entered table) and leave the selection on the second row, second column,
# matching the saved view state of the new sheet.
$newSheet.Range("A1:F2").EntireColumn.AutoFit() | Out-Null
$newSheet.Range("B2").Select() | Out-Null

Write-Host "Added NewContact sheet with header + contact row"
